# LuxMeter blueprints/components.xlsx update:
#   - add "Raspberry Pi Pico 2 W" + a second "SSR Rele Modul 2 kanal" line
#     (two extra corridor-box components), plus two "odpor" (resistor) notes
#   - update selection / scroll position to reflect where the user ended up

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Row 45: second "SSR Rele Modul 2 kanal 5VDC Low Level" line ----------
# Reuse row 19's formatting (same layout: currency C, centered D, bold-
# currency formula E, hyperlink-styled F) so the new cells pick up the
# existing style indices instead of inventing new ones.
$ws.Range("C19:F19").Copy() | Out-Null
$ws.Range("C45:F45").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("B45").Value = "SSR Relé Modul 2 kanál 5VDC Low Level"
$ws.Range("C45").Value = 76
$ws.Range("D45").Value = 3
$ws.Range("E45").Formula = "=C45*D45"
$ws.Range("F45").Value = "https://dratek.cz/arduino/1347-ssr-rele-modul-2-kanaly-5vdc-250vac-omron-g3mb-202p-solid-state-pro-arduino.html"
$ws.Range("F45").Hyperlinks.Add($ws.Range("F45"), "https://dratek.cz/arduino/1347-ssr-rele-modul-2-kanaly-5vdc-250vac-omron-g3mb-202p-solid-state-pro-arduino.html") | Out-Null
# Adding the hyperlink re-applied the built-in "Hyperlink" cell style; copy
# F19's formatting back over so F45 keeps the same look as the rest of the
# sheet's link cells.
$ws.Range("F19").Copy() | Out-Null
$ws.Range("F45").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

# --- Row 46: new "Raspberry Pi Pico 2 W" line ------------------------------
$ws.Range("C19:F19").Copy() | Out-Null
$ws.Range("C46:F46").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

# Set the link text/URL before the product name so new shared strings are
# registered URL-first, matching the source order.
$ws.Range("F46").Value = "https://rpishop.cz/554053/raspberry-pi-pico-2-w "
$ws.Range("B46").Value = "Raspberry Pi Pico 2 W"
$ws.Range("C46").Value = 209
$ws.Range("D46").Value = 10
$ws.Range("E46").Formula = "=C46*D46"
$ws.Range("F46").Hyperlinks.Add($ws.Range("F46"), "https://rpishop.cz/554053/raspberry-pi-pico-2-w ") | Out-Null
$ws.Range("F19").Copy() | Out-Null
$ws.Range("F46").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

# --- Rows 47-48: resistor notes for the new Pico 2 W wiring ---------------
$ws.Range("B47").Value = "odpor 33k"
$ws.Range("B48").Value = "odpor 10k"

# --- Reflect where the user ended up: scrolled down, C48 selected ---------
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C48").Select() | Out-Null
